$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Sheet1"

$srcRange = $ws1.Range("A9:K40")
$srcRange.Copy()
$dstRange = $newSheet.Range("A1")
$dstRange.PasteSpecial("All")
Write-Output "done"
